# Adds a "2022-Q4" quarter sheet to the workbook:
#   - inserts a new summary row on the "总计" sheet for 2022-Q4
#   - duplicates the old "2022-Q2" sheet (so its historical data is preserved
#     under the "2022-Q2" name) and turns the original sheet into the new
#     "2022-Q4" sheet with this quarter's fund holdings

$wb = $excel.ActiveWorkbook

# Helper: force a value onto a range as literal TEXT (not auto-converted to a
# number), without minting a new cell style (no quote-prefix / NumberFormat
# changes). It works by writing a string-literal formula into a scratch cell,
# copying its computed (text) result, and pasting values-only onto the
# destination - this mirrors how Excel turns a formula result into a plain
# value, and keeps the destination on the shared-string/text type.
function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    $scratch = $Range.Worksheet.Range("ZZ100")
    $scratch.Formula = '="' + $Text + '"'
    $scratch.Copy()
    $Range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for "2022-Q4", pushing the existing
#    quarters down by one row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Range("A5:D5").Copy($wsTotal.Range("A6:D6"))
$wsTotal.Range("A4:D4").Copy($wsTotal.Range("A5:D5"))
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4:D4"))
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.06

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q2" sheet right after itself, so the
#    duplicate keeps the old quarter's data, then repurpose the original
#    sheet (still in 2nd tab position) as the brand-new "2022-Q4" sheet.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wb.Worksheets.Item(3))
$wsQ2Copy = $wb.Worksheets.Item(3)

$wsQ2.Name = "2022-Q4temp"
$wsQ2Copy.Name = "2022-Q2"
$wsQ2.Name = "2022-Q4"

$wsQ4 = $wsQ2

# ---------------------------------------------------------------------------
# 3. Replace the "2022-Q4" sheet's fund rows with this quarter's holdings.
#    A2/A3 (0/1) and H2/H3 (6/6) are unchanged; B-G change to the new fund.
# ---------------------------------------------------------------------------
Set-TextValue $wsQ4.Range("B2") "008905"
Set-TextValue $wsQ4.Range("C2") "嘉合锦鹏添利混合A"
Set-TextValue $wsQ4.Range("D2") "3.64"
Set-TextValue $wsQ4.Range("E2") "22.86"
Set-TextValue $wsQ4.Range("F2") "0.97"
Set-TextValue $wsQ4.Range("G2") "0.0353"

Set-TextValue $wsQ4.Range("B3") "008906"
Set-TextValue $wsQ4.Range("C3") "嘉合锦鹏添利混合C"
Set-TextValue $wsQ4.Range("D3") "2.16"
Set-TextValue $wsQ4.Range("E3") "22.86"
Set-TextValue $wsQ4.Range("F3") "0.97"
Set-TextValue $wsQ4.Range("G3") "0.0210"

# ---------------------------------------------------------------------------
# 4. Restore the active-sheet/tab-selected state to the last sheet, matching
#    the original workbook (avoids leaving the copy step's activation as a
#    stray diff).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
